# excel data source files updates
# - add tyre_before / tyre_after columns (H/I) to the "pitstop" sheet
# - add a new "weather" sheet with race-day weather data
# - add a new "altitude" sheet with a single delta value

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. pitstop: add tyre_before (H) / tyre_after (I) columns
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("pitstop")
[void]$ws4.Activate()

$ws4.Range("H1").Value = "tyre_before"
$ws4.Range("I1").Value = "tyre_after"

$hiData = @(
  @(3,2),
  @(2,1),
  @(3,2),
  @(2,1),
  @(3,1),
  @(1,2),
  @(2,1),
  @(3,1),
  @(2,1),
  @(2,1),
  @(1,3),
  @(3,3),
  @(3,2),
  @(2,1),
  @(2,1),
  @(1,3),
  @(3,1),
  @(3,2),
  @(3,2),
  @(2,1),
  @(3,1),
  @(2,2),
  @(2,1),
  @(2,1),
  @(2,1),
  @(2,1),
  @(2,2),
  @(2,3),
  @(2,1)
)

for ($idx = 0; $idx -lt $hiData.Count; $idx++) {
    $row = $idx + 2
    $ws4.Range("H$row").Value = $hiData[$idx][0]
    $ws4.Range("I$row").Value = $hiData[$idx][1]
}

$ws4.Columns.Item(8).ColumnWidth = 10.666666666666666
$ws4.Columns.Item(9).ColumnWidth = 9.0

[void]$ws4.Range("H31").Select()

# ---------------------------------------------------------------------------
# 2. weather: new sheet with sky/temperature/humidity/wind data
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws8.Name = "weather"

$ws8.Range("A1").Value = "Skycondition"
$ws8.Range("B1").Value = "Mostly Cloudy"
$ws8.Range("A2").Value = "Temperature"
$ws8.Range("B2").Value = "64.91°F"
$ws8.Range("A3").Value = "Humidity"
$ws8.Range("B3").Value = 0.72
$ws8.Range("B3").NumberFormat = "0%"
$ws8.Range("A4").Value = "Wind speed"
$ws8.Range("B4").Value = "6.3 mph"
$ws8.Range("A5").Value = "Wind bearing"
$ws8.Range("B5").Value = "27°"

$ws8.Columns.Item(1).ColumnWidth = 12.0
$ws8.Columns.Item(2).ColumnWidth = 12.833333333333334

$ws8.PageSetup.Orientation = 1

[void]$ws8.Range("F2").Select()

# ---------------------------------------------------------------------------
# 3. altitude: new sheet with a single delta value
# ---------------------------------------------------------------------------
$ws9 = $wb.Worksheets.Add([System.Type]::Missing, $ws8)
$ws9.Name = "altitude"

$ws9.Range("A1").Value = "delta"
$ws9.Range("B1").Value = 11.3

[void]$ws9.Range("B2").Select()
[void]$ws9.Activate()
